# Improve logging system configuration: append the latest captured log
# record (row 60) to each of the MID database worksheets.

$wb = $excel.ActiveWorkbook

$rows = @{
    "MID_LFT_#1" = @{
        A = 45846.46619212963
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x60"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 352
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45846.46619212963
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x5C"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 348
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45846.46619212963
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x68"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 104
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45846.46619212963
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7D"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 125
        I = 9
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rows[$ws.Name]
    if ($data -ne $null) {
        $newRow = 60

        $ws.Cells.Item($newRow, 1).Value = $data.A
        $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $ws.Cells.Item($newRow, 2).Value = $data.B
        $ws.Cells.Item($newRow, 3).Value = $data.C
        $ws.Cells.Item($newRow, 4).Value = $data.D
        $ws.Cells.Item($newRow, 5).Value = $data.E

        $ws.Cells.Item($newRow, 6).Value = $data.F
        $ws.Cells.Item($newRow, 7).Value = [double]$data.G
        $ws.Cells.Item($newRow, 8).Value = $data.H
        $ws.Cells.Item($newRow, 9).Value = $data.I
    }
}
